# Generate Report for Handback
# Adds a new handback record (f314184d-8bc2-4f71-904f-842bebebe07c) as row 3
# on each sheet, and updates the existing record's uuid/timestamps from
# 9370d91d-79c1-4bf4-94db-e856dd5822b6 to cfbe7b37-69a7-442d-8129-af623d0036d4.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$oldUuid = "9370d91d-79c1-4bf4-94db-e856dd5822b6"
$newUuid1 = "cfbe7b37-69a7-442d-8129-af623d0036d4"
$newUuid2 = "f314184d-8bc2-4f71-904f-842bebebe07c"

$mainRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f21bce2d6b9fdaec93f76fe0ea12853afc607609/e2e/"
$zhcnRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8888789d2b14b2c0be250732d975e5d2ad6fbf29/e2e/"
$dedeRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/38a0d9643237675fd8c07dcb8ca602e780ff7184/e2e/"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# NOTE: Range.Hyperlinks.Delete() on this runtime clears *all* hyperlinks on
# the worksheet (not just the target range), so we delete once up front and
# re-Add every hyperlink (existing + new) afterwards, in document order.
$wsOv.Range("B2").Hyperlinks.Delete() | Out-Null

# Row 2: refresh uuid + timestamp on the existing record
$wsOv.Range("A2").Value2 = "$newUuid1.md"
$wsOv.Range("G2").NumberFormat = $dateFmt
$wsOv.Range("G2").Value2 = "2016-08-16 13:01:26"

# Row 3: new record
$wsOv.Range("A3").Value2 = "$newUuid2.md"
$wsOv.Range("C3").Value2 = ".md"
$wsOv.Range("E3").Value2 = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value2 = "Handed back: in sync with en-US"
$wsOv.Range("G3").NumberFormat = $dateFmt
$wsOv.Range("G3").Value2 = "2016-08-16 13:01:26"

# Re-create hyperlinks in document order: B2 (existing, refreshed), B3 (new)
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "$mainRepoBase$newUuid1.md", "", "", "e2e\$newUuid1.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "$mainRepoBase$newUuid2.md", "", "", "e2e\$newUuid2.md") | Out-Null

$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Clear all hyperlinks up front (Range.Hyperlinks.Delete() clears the whole
# sheet on this runtime), they are all re-Added below in document order.
$wsZh.Range("A2").Hyperlinks.Delete() | Out-Null

# Row 2: refresh uuid + hashes + timestamps on the existing record
$wsZh.Range("G2").Value2 = "$newUuid1.91245b5c9e08bc27eaa927ccf13be4fbe9e3b8f7.zh-cn.xlf"
$wsZh.Range("H2").NumberFormat = $dateFmt
$wsZh.Range("H2").Value2 = "2016-08-16 13:01:20"
$wsZh.Range("J2").Value2 = "$newUuid1.91245b5c9e08bc27eaa927ccf13be4fbe9e3b8f7.zh-cn.xlf"
$wsZh.Range("K2").NumberFormat = $dateFmt
$wsZh.Range("K2").Value2 = "2016-08-16 13:01:48"

# Row 3: new record
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value2 = "e2e"
$wsZh.Range("E3").Value2 = "ht"
$wsZh.Range("F3").Value2 = "'True"
$wsZh.Range("G3").Value2 = "$newUuid2.cb5b91a0a1ac2e0ad5ba66308a0eceb605ec02b2.zh-cn.xlf"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("H3").Value2 = "2016-08-16 13:01:20"
$wsZh.Range("J3").Value2 = "$newUuid2.cb5b91a0a1ac2e0ad5ba66308a0eceb605ec02b2.zh-cn.xlf"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("K3").Value2 = "2016-08-16 13:01:48"
$wsZh.Range("L3").Value2 = "'"
$wsZh.Range("M3").Value2 = "'True"
$wsZh.Range("N3").Value2 = "'"
$wsZh.Range("O3").Value2 = "'False"
$wsZh.Range("P3").Value2 = "'"

# Re-create hyperlinks in document order: A2, I2 (existing, refreshed), A3, I3 (new)
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$mainRepoBase$newUuid1.md", "", "", "$newUuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$zhcnRepoBase$newUuid1.md", "", "", "$newUuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$mainRepoBase$newUuid2.md", "", "", "$newUuid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$zhcnRepoBase$newUuid2.md", "", "", "$newUuid2.md") | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Clear all hyperlinks up front (Range.Hyperlinks.Delete() clears the whole
# sheet on this runtime), they are all re-Added below in document order.
$wsDe.Range("A2").Hyperlinks.Delete() | Out-Null

# Row 2: refresh uuid + hashes + timestamps on the existing record
$wsDe.Range("G2").Value2 = "$newUuid1.91245b5c9e08bc27eaa927ccf13be4fbe9e3b8f7.de-de.xlf"
$wsDe.Range("H2").NumberFormat = $dateFmt
$wsDe.Range("H2").Value2 = "2016-08-16 13:01:26"
$wsDe.Range("J2").Value2 = "$newUuid1.91245b5c9e08bc27eaa927ccf13be4fbe9e3b8f7.de-de.xlf"
$wsDe.Range("K2").NumberFormat = $dateFmt
$wsDe.Range("K2").Value2 = "2016-08-16 13:01:56"

# Row 3: new record
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value2 = "e2e"
$wsDe.Range("E3").Value2 = "ht"
$wsDe.Range("F3").Value2 = "'True"
$wsDe.Range("G3").Value2 = "$newUuid2.cb5b91a0a1ac2e0ad5ba66308a0eceb605ec02b2.de-de.xlf"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("H3").Value2 = "2016-08-16 13:01:26"
$wsDe.Range("J3").Value2 = "$newUuid2.cb5b91a0a1ac2e0ad5ba66308a0eceb605ec02b2.de-de.xlf"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("K3").Value2 = "2016-08-16 13:01:56"
$wsDe.Range("L3").Value2 = "'"
$wsDe.Range("M3").Value2 = "'True"
$wsDe.Range("N3").Value2 = "'"
$wsDe.Range("O3").Value2 = "'False"
$wsDe.Range("P3").Value2 = "'"

# Re-create hyperlinks in document order: A2, I2 (existing, refreshed), A3, I3 (new)
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$mainRepoBase$newUuid1.md", "", "", "$newUuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$dedeRepoBase$newUuid1.md", "", "", "$newUuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$mainRepoBase$newUuid2.md", "", "", "$newUuid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$dedeRepoBase$newUuid2.md", "", "", "$newUuid2.md") | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3")) | Out-Null

Write-Host "Handback report rows generated."
